$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @("Price (D)", "Volume 1h (E)")
$updates = @{
    2 = @("27.934.03", "  +0.70%  ")
    3 = @("1.880.44", "  +0.08%  ")
    4 = @("1.018", "  +1.56%  ")
    5 = @("335.20", "  +1.15%  ")
    6 = @("1.017", "  +1.44%  ")
    7 = @("0.4708", "  -0.07%  ")
    8 = @("0.3918", "  -1.09%  ")
    9 = @("46.76", "  -2.26%  ")
    10 = @("0.07957", "  -0.92%  ")
    11 = @("1.008", "  -1.48%  ")
    12 = @("21.66", "  -1.13%  ")
    13 = @("1.867.95", "  -1.49%  ")
    14 = @("5.952", "  -0.30%  ")
    15 = @("7.103", "  -0.79%  ")
    16 = @("1.020", "  +1.47%  ")
    17 = @("0.06784", "  +2.66%  ")
    18 = @("87.40", "  +0.30%  ")
    19 = @("0.00001046", "  +0.17%  ")
    20 = @("17.05", "  -0.92%  ")
    21 = @("1.018", "  +1.52%  ")
    22 = @("27.937.94", "  +0.85%  ")
    23 = @("5.474", "  -0.61%  ")
    24 = @("10.94", "  -0.84%  ")
    25 = @("2.360", "  +2.77%  ")
    26 = @("2.109.29", "  -0.23%  ")
    27 = @("160.05", "  +2.25%  ")
    28 = @("19.91", "  -1.61%  ")
    29 = @("2.080", "  -0.64%  ")
    30 = @("5.458", "  -2.43%  ")
    31 = @("121.01", "  -1.30%  ")
    32 = @("0.09531", "  -0.25%  ")
    33 = @("0.9604", "  -1.25%  ")
    34 = @("3.659", "  +0.85%  ")
    35 = @("5.319", "  +0.23%  ")
    36 = @("1.346", "  -7.50%  ")
    37 = @("0.06112", "  -0.13%  ")
    38 = @("0.02240", "  -1.13%  ")
    39 = @("1.202", "  -2.68%  ")
    40 = @("1.016", "  +1.39%  ")
    41 = @("8.121", "  -0.55%  ")
    42 = @("0.5913", "  -1.48%  ")
    43 = @("0.1891", "  -0.47%  ")
    44 = @("10.19", "  -0.35%  ")
    45 = @("1.269", "  +1.53%  ")
    46 = @("0.5651", "  -0.61%  ")
    47 = @("12.13", "  -0.78%  ")
    48 = @("3.396", "  -0.12%  ")
    49 = @("1.919", "  -0.74%  ")
    50 = @("0.06866", "  +0.65%  ")
    51 = @("113.74", "  +1.88%  ")
}

foreach ($row in $updates.Keys) {
    $priceText = $updates[$row][0]
    $volumeText = $updates[$row][1]

    $dCell = $ws.Range("D" + $row)
    # Force the "Price" cell to stay text (avoid Excel auto-converting
    # numeric-looking strings like "1.018" into a real number), then
    # strip the temporary number format so no extra style is introduced.
    $dCell.NumberFormat = "@"
    $dCell.Value = $priceText
    $dCell.ClearFormats()

    $ws.Range("E" + $row).Value = $volumeText
}
